$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Notes:"
$ws.Range("B2").Value = 26
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 11
$ws.Range("F2").Value = 11
$ws.Range("G2").Value = "Run on the easy board, average of ten trials"

$ws.Columns.Item(7).ColumnWidth = 55.16

$ws.Range("B3").Select()
